# Insert two new rows ("on"/"ON" and "off"/"OFF" placeholder options) right
# after the header-ish rows, before the existing "cycle" row — pushing the
# rest of the language table down by two rows (quiz modal placeholder
# options).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("6:7").Insert() | Out-Null

$ws.Range("A6").Value2 = "on"
$ws.Range("B6").Value2 = "ON"

$ws.Range("A7").Value2 = "off"
$ws.Range("B7").Value2 = "OFF"

# Update selection / scroll position to match the authored state.
$ws.Range("A7").Select() | Out-Null
